$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "IP" header (H1) onto the two new
# header cells so they pick up the same bold/border/center style (s="1").
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).Value = 6
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 5
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(4, 10).Value = 6
$ws.Cells.Item(5, 9).Value = 1
$ws.Cells.Item(5, 10).Value = 6
$ws.Cells.Item(6, 9).Value = 1
$ws.Cells.Item(6, 10).Value = 3
$ws.Cells.Item(7, 9).Value = 1
$ws.Cells.Item(7, 10).Value = 6
$ws.Cells.Item(8, 9).Value = 1
$ws.Cells.Item(8, 10).Value = 6
$ws.Cells.Item(9, 9).Value = 1
$ws.Cells.Item(9, 10).Value = 6
$ws.Cells.Item(10, 9).Value = 1
$ws.Cells.Item(10, 10).Value = 6
$ws.Cells.Item(11, 9).Value = 1
$ws.Cells.Item(11, 10).Value = 5
$ws.Cells.Item(12, 9).Value = 1
$ws.Cells.Item(12, 10).Value = 5
$ws.Cells.Item(13, 9).Value = 1
$ws.Cells.Item(13, 10).Value = 4
$ws.Cells.Item(14, 9).Value = 8
$ws.Cells.Item(14, 10).Value = 9
$ws.Cells.Item(15, 9).Value = 8
$ws.Cells.Item(15, 10).Value = 8
$ws.Cells.Item(16, 9).Value = 7
$ws.Cells.Item(16, 10).Value = 8
$ws.Cells.Item(17, 9).Value = 7
$ws.Cells.Item(17, 10).Value = 8
$ws.Cells.Item(18, 9).Value = 6
$ws.Cells.Item(18, 10).Value = 7
$ws.Cells.Item(19, 9).Value = 8
$ws.Cells.Item(19, 10).Value = 9
$ws.Cells.Item(20, 9).Value = 8
$ws.Cells.Item(20, 10).Value = 9
$ws.Cells.Item(21, 9).Value = 9
$ws.Cells.Item(21, 10).Value = 9
$ws.Cells.Item(22, 9).Value = 6
$ws.Cells.Item(22, 10).Value = 8
$ws.Cells.Item(23, 9).Value = 9
$ws.Cells.Item(23, 10).Value = 9
$ws.Cells.Item(24, 9).Value = 8
$ws.Cells.Item(24, 10).Value = 8
$ws.Cells.Item(25, 9).Value = 6
$ws.Cells.Item(25, 10).Value = 7
$ws.Cells.Item(26, 9).Value = 1
$ws.Cells.Item(26, 10).Value = 4
$ws.Cells.Item(27, 9).Value = 1
$ws.Cells.Item(27, 10).Value = 6
$ws.Cells.Item(28, 9).Value = 1
$ws.Cells.Item(28, 10).Value = 3
$ws.Cells.Item(29, 9).Value = 1
$ws.Cells.Item(29, 10).Value = 5
$ws.Cells.Item(30, 9).Value = 6
$ws.Cells.Item(30, 10).Value = 9
$ws.Cells.Item(31, 9).Value = 1
$ws.Cells.Item(31, 10).Value = 4
$ws.Cells.Item(32, 9).Value = 1
$ws.Cells.Item(32, 10).Value = 5
$ws.Cells.Item(33, 9).Value = 1
$ws.Cells.Item(33, 10).Value = 4
$ws.Cells.Item(34, 9).Value = 1
$ws.Cells.Item(34, 10).Value = 3
$ws.Cells.Item(35, 9).Value = 1
$ws.Cells.Item(35, 10).Value = 2
